{"js": "const replacements = [\n  [\"2024-06-29 Saturday\", \"2024-06-30 Sunday\"],\n  [\"696\u00d72=1392\", \"726\u00d75=3630\"],\n  [\"999\u00d78=7992\", \"511\u00d75=2555\"],\n  [\"367\u00d75=1835\", \"611\u00d74=2444\"],\n  [\"272\u00d73=816\", \"499\u00d79=4491\"],\n  [\"452\u00d76=2712\", \"262\u00d72=524\"],\n  [\"414\u00d72=828\", \"199\u00d74=796\"],\n  [\"878\u00d76=5268\", \"166\u00d72=332\"],\n  [\"714\u00d72=1428\", \"980\u00d76=5880\"],\n  [\"979\u00d76=5874\", \"231\u00d76=1386\"],\n  [\"537\u00d77=3759\", \"704\u00d74=2816\"],\n  [\"982\u00d73=2946\", \"798\u00d78=6384\"],\n  [\"246\u00d78=1968\", \"471\u00d79=4239\"],\n  [\"744\u00d79=6696\", \"876\u00d72=1752\"],\n  [\"215\u00d72=430\", \"511\u00d76=3066\"],\n  [\"727\u00d79=6543\", \"713\u00d76=4278\"],\n  [\"248\u00d77=1736\", \"738\u00d79=6642\"],\n  [\"998\u00d77=6986\", \"825\u00d72=1650\"],\n  [\"193\u00d77=1351\", \"317\u00d74=1268\"],\n  [\"477\u00d75=2385\", \"888\u00d78=7104\"],\n  [\"621\u00d78=4968\", \"372\u00d77=2604\"],\n  [\"864\u00d73=2592\", \"151\u00d76=906\"],\n  [\"136\u00d73=408\", \"812\u00d73=2436\"],\n  [\"781\u00d77=5467\", \"858\u00d76=5148\"],\n  [\"922\u00d74=3688\", \"621\u00d76=3726\"],\n  [\"737\u00d77=5159\", \"415\u00d72=830\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-06-29 Saturday\", \"2024-06-30 Sunday\"),\n    @(\"696\u00d72=1392\", \"726\u00d75=3630\"),\n    @(\"999\u00d78=7992\", \"511\u00d75=2555\"),\n    @(\"367\u00d75=1835\", \"611\u00d74=2444\"),\n    @(\"272\u00d73=816\", \"499\u00d79=4491\"),\n    @(\"452\u00d76=2712\", \"262\u00d72=524\"),\n    @(\"414\u00d72=828\", \"199\u00d74=796\"),\n    @(\"878\u00d76=5268\", \"166\u00d72=332\"),\n    @(\"714\u00d72=1428\", \"980\u00d76=5880\"),\n    @(\"979\u00d76=5874\", \"231\u00d76=1386\"),\n    @(\"537\u00d77=3759\", \"704\u00d74=2816\"),\n    @(\"982\u00d73=2946\", \"798\u00d78=6384\"),\n    @(\"246\u00d78=1968\", \"471\u00d79=4239\"),\n    @(\"744\u00d79=6696\", \"876\u00d72=1752\"),\n    @(\"215\u00d72=430\", \"511\u00d76=3066\"),\n    @(\"727\u00d79=6543\", \"713\u00d76=4278\"),\n    @(\"248\u00d77=1736\", \"738\u00d79=6642\"),\n    @(\"998\u00d77=6986\", \"825\u00d72=1650\"),\n    @(\"193\u00d77=1351\", \"317\u00d74=1268\"),\n    @(\"477\u00d75=2385\", \"888\u00d78=7104\"),\n    @(\"621\u00d78=4968\", \"372\u00d77=2604\"),\n    @(\"864\u00d73=2592\", \"151\u00d76=906\"),\n    @(\"136\u00d73=408\", \"812\u00d73=2436\"),\n    @(\"781\u00d77=5467\", \"858\u00d76=5148\"),\n    @(\"922\u00d74=3688\", \"621\u00d76=3726\"),\n    @(\"737\u00d77=5159\", \"415\u00d72=830\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
